# Update countries & provincias Spain
# This refreshes the COVID-19 country case counts on the "Pais" sheet and
# re-ranks a handful of countries whose total-cases ordering changed as a
# result (their rows keep their position, but the country name + stats
# for that position shift, same as the source data export does).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 11:22"

# Estados Unidos
$ws.Range("B4").Value = 2330769
$ws.Range("C4").Value = 191
$ws.Range("D4").Value = 973055
$ws.Range("E4").Value = 1235731
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 121983

# India
$ws.Range("B7").Value = 412788
$ws.Range("C7").Value = 1061
$ws.Range("D7").Value = 228504
$ws.Range("E7").Value = 170994
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 13
$ws.Range("H7").Value = 13290

# Banglades
$ws.Range("B20").Value = 112306
$ws.Range("C20").Value = 3531
$ws.Range("D20").Value = 45077
$ws.Range("E20").Value = 65765
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 39
$ws.Range("H20").Value = 1464

# Bielorrusia
$ws.Range("B32").Value = 45891
$ws.Range("C32").Value = 862
$ws.Range("D32").Value = 18404
$ws.Range("E32").Value = 25022
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 36
$ws.Range("H32").Value = 2465

# Singapur
$ws.Range("B39").Value = 31931
$ws.Range("C39").Value = 311
$ws.Range("D39").Value = 16683
$ws.Range("E39").Value = 13892
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 10
$ws.Range("H39").Value = 1356

# Rows 41-44 re-rank: Oman overtakes Filipinas, pushing
# Filipinas / Irak / Afganistan down one position each.
$ws.Range("A41").Value = "Oman"
$ws.Range("B41").Value = 29471
$ws.Range("C41").Value = 905
$ws.Range("D41").Value = 15552
$ws.Range("E41").Value = 13788
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 131

$ws.Range("A42").Value = "Filipinas"
$ws.Range("B42").Value = 29400
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 7650
$ws.Range("E42").Value = 20600
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 1150

$ws.Range("A43").Value = "Irak"
$ws.Range("B43").Value = 29222
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 13211
$ws.Range("E43").Value = 14998
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 1013

$ws.Range("A44").Value = "Afganistan"
$ws.Range("B44").Value = 28833
$ws.Range("C44").Value = 409
$ws.Range("D44").Value = 8764
$ws.Range("E44").Value = 19488
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 12
$ws.Range("H44").Value = 581

# Bolivia
$ws.Range("B50").Value = 21331
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 15790
$ws.Range("E50").Value = 5480
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 61

# Israel
$ws.Range("B51").Value = 20686
$ws.Range("C51").Value = 53
$ws.Range("D51").Value = 15664
$ws.Range("E51").Value = 4717
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 305

# Guatemala
$ws.Range("B56").Value = 17225
$ws.Range("C56").Value = 446
$ws.Range("D56").Value = 10719
$ws.Range("E56").Value = 6388
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 118

# Costa de Marfil
$ws.Range("B72").Value = 8572
$ws.Range("C72").Value = 16
$ws.Range("D72").Value = 8156
$ws.Range("E72").Value = 295
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 121

# Tayikistan
$ws.Range("B75").Value = 7143
$ws.Range("C75").Value = 1
$ws.Range("D75").Value = 6200
$ws.Range("E75").Value = 617
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 326

# Rows 97-98 re-rank: Republica de Africa Central overtakes Somalia.
$ws.Range("A97").Value = "Republica de Africa Central"
$ws.Range("B97").Value = 2808
$ws.Range("C97").Value = 122
$ws.Range("D97").Value = 472
$ws.Range("E97").Value = 2313
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 4
$ws.Range("H97").Value = 23

$ws.Range("A98").Value = "Somalia"
$ws.Range("B98").Value = 2755
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 751
$ws.Range("E98").Value = 1916
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 88

# Republica de Yibuti
$ws.Range("B112").Value = 1798
$ws.Range("C112").Value = 3
$ws.Range("D112").Value = 1475
$ws.Range("E112").Value = 247
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 76

# Swazilandia
$ws.Range("B137").Value = 827
$ws.Range("C137").Value = 43
$ws.Range("D137").Value = 439
$ws.Range("E137").Value = 385
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 3

# Liberia
$ws.Range("B138").Value = 770
$ws.Range("C138").Value = 7
$ws.Range("D138").Value = 578
$ws.Range("E138").Value = 192
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0

# Row 187 country
$ws.Range("B187").Value = 55
$ws.Range("C187").Value = 9
$ws.Range("D187").Value = 19
$ws.Range("E187").Value = 36
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0

# Rows 202-203 re-rank: Dominica overtakes Fiyi (same stats, only the
# country names swap places).
$ws.Range("A202").Value = "Dominica"
$ws.Range("B202").Value = 18
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 18
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("A203").Value = "Fiyi"
$ws.Range("B203").Value = 18
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 18
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Rows 208-209 re-rank: Islas Turcas y Caicos overtakes Santa Sede.
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("B208").Value = 12
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 11
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("B209").Value = 12
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 12
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 7
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 8
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
